$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Widen column B (42.18359375 -> 52.0 in the saved XML "width" attribute).
#    The ColumnWidth COM property is offset by 5/6 (0.8333...) from the raw
#    OOXML column width, so back that out to land exactly on 52.0.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 52.0 - 0.8333333333333333

# ---------------------------------------------------------------------------
# 2. Insert 4 new rows before the current row 4 (RULE-100), pushing the
#    existing rows 4-8 down to rows 8-12. The new rows inherit the style
#    (s="2", the PASS/green look) of row 3 above them, same as the target.
# ---------------------------------------------------------------------------
$ws.Rows("4:7").Insert()

# New row 4: RULE-107
$ws.Range("A4").Value = "RULE-107"
$ws.Range("B4").Value = "Forbidden substring check for .properties files"
$ws.Range("C4").Value = "HIGH"
$ws.Range("D4").Value = "PASS"
$ws.Range("E4").Value = "All checks passed"

# New row 5: RULE-108
$ws.Range("A5").Value = "RULE-108"
$ws.Range("B5").Value = "Forbidden substring check for .policy files"
$ws.Range("C5").Value = "HIGH"
$ws.Range("D5").Value = "PASS"
$ws.Range("E5").Value = "All checks passed"

# New row 6: RULE-109
$ws.Range("A6").Value = "RULE-109"
$ws.Range("B6").Value = "Forbidden regex pattern (ip addresses) check in .properties files"
$ws.Range("C6").Value = "HIGH"
$ws.Range("D6").Value = "PASS"
$ws.Range("E6").Value = "All checks passed"

# New row 7: RULE-110
$ws.Range("A7").Value = "RULE-110"
$ws.Range("B7").Value = "Forbidden regex pattern (ip addresses) check in .policy files"
$ws.Range("C7").Value = "HIGH"
$ws.Range("D7").Value = "PASS"
$ws.Range("E7").Value = "All checks passed"

# ---------------------------------------------------------------------------
# 3. Wording tweaks in the detail columns of the rows that used to be 4-7
#    (now rows 8-11 after the insert above).
# ---------------------------------------------------------------------------

# Row 8 (was row 4): RULE-100 - "Token" -> "Required token"
$ws.Range("E8").Value = "• Validation failures:`n• Required token 'apiId' not found in file: Properties\OCP\ITE.properties (case-sensitive: true)"
$ws.Range("E8").EntireRow.AutoFit()

# Row 9 (was row 5): RULE-101 - drop ", test" from the expected-values list
$ws.Range("E9").Value = "• Validation failures:`n• Property 'LogJsonFormat' found but value does not match expected values [true, false] in file: Properties\OCP\ITE.properties`n• Property 'anotherpropertycheck' found but value does not match expected values [somevalue] in file: Properties\OCP\ITE.properties"
$ws.Range("E9").EntireRow.AutoFit()

# Row 10 (was row 6): RULE-102 - "Token" -> "Required token" (x4)
$ws.Range("E10").Value = "• Validation failures:`n• Required token 'http.protocols=HTTPS' not found in file: Policies\TDV.policy (case-sensitive: true)`n• Required token 'http.private.port=8081' not found in file: Policies\TDV.policy (case-sensitive: true)`n• Required token 'http.protocols=HTTPS' not found in file: Policies\TDV1.policy (case-sensitive: true)`n• Required token 'http.private.port=8081' not found in file: Policies\TDV1.policy (case-sensitive: true)"
$ws.Range("E10").EntireRow.AutoFit()

# Rows 11 (RULE-103) and 12 (RULE-104) (were rows 7 and 8) are unchanged in content.
